$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Stash the "last row" border formatting (row 22) in the currently blank
#     row 1 so we can re-apply it once the table shrinks and a different row
#     becomes the new last row.
$ws.Range("B22:J22").Copy()
$ws.Range("B1:J1").PasteSpecial(-4122)

# --- Remove the two obsolete detail rows so the table shrinks from 7
#     employee rows (16-22) down to 5 (16-20). Deleting the same row index
#     twice removes the two trailing rows and shifts everything below (the
#     footer block) up by two rows, matching the target layout.
$ws.Rows.Item(21).Delete()
$ws.Rows.Item(21).Delete()

# --- Re-apply the stashed "last row" formatting onto the new final row (20)
$ws.Range("B1:J1").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)
$ws.Range("B1:J1").Clear()

# --- Refresh the summary figures at the top of the statement
$ws.Range("E11").Value = 160108
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 3

# --- Rewrite the remaining 4 detail rows with the updated account data
# Row 17: Andres Ricardo Perez Agamez, periodo 1805
$ws.Cells.Item(17, 3).Value = "1143369763"
$ws.Cells.Item(17, 4).Value = "ANDRES RICARDO PEREZ AGAMEZ"
$ws.Cells.Item(17, 5).Value = "1805"
$ws.Cells.Item(17, 6).Value = 31249
$ws.Cells.Item(17, 7).Value = 877803

# Row 18: Ana Rita de la Rosa Rondon, periodo 1901
$ws.Cells.Item(18, 3).Value = "45500352"
$ws.Cells.Item(18, 4).Value = "ANA RITA DE LA ROSA RONDON"
$ws.Cells.Item(18, 5).Value = "1901"
$ws.Cells.Item(18, 6).Value = 31249
$ws.Cells.Item(18, 7).Value = 781242

# Row 19: Duglas Robles Niño, periodo 1901
$ws.Cells.Item(19, 3).Value = "1052988508"
$ws.Cells.Item(19, 4).Value = "DUGLAS ROBLES NIÑO"
$ws.Cells.Item(19, 5).Value = "1901"
$ws.Cells.Item(19, 6).Value = 31249
$ws.Cells.Item(19, 7).Value = 781242

# Row 20: Luisa Fernanda Villalba Padilla, periodo 2005
$ws.Cells.Item(20, 3).Value = "1050958772"
$ws.Cells.Item(20, 4).Value = "LUISA FERNANDA VILLALBA PADILLA"
$ws.Cells.Item(20, 5).Value = "2005"
$ws.Cells.Item(20, 6).Value = 35112
$ws.Cells.Item(20, 7).Value = 877803
